$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 2 is updated, and new rows 3-7 are appended, matching the
# refreshed NATMI ligand-receptor pairing table (Spn -> Siglec1) across the
# ECs / M2 sending & target clusters.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Spn"
$ws.Range("C2").Value = "Siglec1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7483270000000001
$ws.Range("H2").Value = 2.244981
$ws.Range("I2").Value = 0.2123245966460564
$ws.Range("J2").Value = 0.2123245966460564
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 4.511004
$ws.Range("N2").Value = 13.533012
$ws.Range("O2").Value = 0.05275259630406119
$ws.Range("P2").Value = 0.05275259630406118
$ws.Range("Q2").Value = 3.375706090308
$ws.Range("R2").Value = 30.381354812772
$ws.Range("S2").Value = 0.01120067373229204
$ws.Range("T2").Value = 0.01120067373229204

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Spn"
$ws.Range("C3").Value = "Siglec1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.7483270000000001
$ws.Range("H3").Value = 2.244981
$ws.Range("I3").Value = 0.2123245966460564
$ws.Range("J3").Value = 0.2123245966460564
$ws.Range("K3").Value = 1.0
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.03801766666666666
$ws.Range("N3").Value = 0.114053
$ws.Range("O3").Value = 0.0004445863098523146
$ws.Range("P3").Value = 0.0004445863098523145
$ws.Range("Q3").Value = 0.02844964644366667
$ws.Range("R3").Value = 0.256046817993
$ws.Range("S3").Value = 0.00009439660891375133
$ws.Range("T3").Value = 0.00009439660891375131

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Spn"
$ws.Range("C4").Value = "Siglec1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.7483270000000001
$ws.Range("H4").Value = 2.244981
$ws.Range("I4").Value = 0.2123245966460564
$ws.Range("J4").Value = 0.2123245966460564
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 80.963433
$ws.Range("N4").Value = 242.890299
$ws.Range("O4").Value = 0.9468028173860865
$ws.Range("P4").Value = 0.9468028173860864
$ws.Range("Q4").Value = 60.587122926591
$ws.Range("R4").Value = 545.284106339319
$ws.Range("S4").Value = 0.2010295263048506
$ws.Range("T4").Value = 0.2010295263048506

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Spn"
$ws.Range("C5").Value = "Siglec1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 2.776121
$ws.Range("H5").Value = 8.328363
$ws.Range("I5").Value = 0.7876754033539436
$ws.Range("J5").Value = 0.7876754033539437
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 4.511004
$ws.Range("N5").Value = 13.533012
$ws.Range("O5").Value = 0.05275259630406119
$ws.Range("P5").Value = 0.05275259630406118
$ws.Range("Q5").Value = 12.523092935484
$ws.Range("R5").Value = 112.707836419356
$ws.Range("S5").Value = 0.04155192257176915
$ws.Range("T5").Value = 0.04155192257176915

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Spn"
$ws.Range("C6").Value = "Siglec1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 2.776121
$ws.Range("H6").Value = 8.328363
$ws.Range("I6").Value = 0.7876754033539436
$ws.Range("J6").Value = 0.7876754033539437
$ws.Range("K6").Value = 1.0
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.03801766666666666
$ws.Range("N6").Value = 0.114053
$ws.Range("O6").Value = 0.0004445863098523146
$ws.Range("P6").Value = 0.0004445863098523145
$ws.Range("Q6").Value = 0.1055416428043333
$ws.Range("R6").Value = 0.9498747852389999
$ws.Range("S6").Value = 0.0003501897009385632
$ws.Range("T6").Value = 0.0003501897009385632

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Spn"
$ws.Range("C7").Value = "Siglec1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 2.776121
$ws.Range("H7").Value = 8.328363
$ws.Range("I7").Value = 0.7876754033539436
$ws.Range("J7").Value = 0.7876754033539437
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 80.963433
$ws.Range("N7").Value = 242.890299
$ws.Range("O7").Value = 0.9468028173860865
$ws.Range("P7").Value = 0.9468028173860864
$ws.Range("Q7").Value = 224.764286583393
$ws.Range("R7").Value = 2022.878579250537
$ws.Range("S7").Value = 0.7457732910812359
$ws.Range("T7").Value = 0.7457732910812359
